# Edit: add a red "(This is a change ...)" suffix to the first paragraph,
# and append a new empty, shaded paragraph after the closing speech paragraph.

$d = $word.ActiveDocument

# --- 1. First paragraph: add two trailing spaces then a red-colored
#        "(This is a change - Version for main branch)" suffix. ---

$firstPara = $d.Paragraphs(1)
$firstRange = $firstPara.Range
# Trim the paragraph mark off the end so we only touch the text run.
$firstRange.MoveEnd(1, -1) | Out-Null
$firstRange.InsertAfter("  ")

# Collapse to an insertion point right after the text + two spaces we just
# added, then insert the red, parenthesized note as new run(s).
$insertPoint = $firstRange.Duplicate
$insertPoint.Collapse(0)

$insertPoint.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$insertPoint.Font.Color = 255

$insertPoint.Collapse(0)
$insertPoint.InsertAfter("rsion for main branch")
$insertPoint.Font.Color = 255

$insertPoint.Collapse(0)
$insertPoint.InsertAfter(")")
$insertPoint.Font.Color = 255

# --- 2. Append a new, empty paragraph with a light-gray shading fill
#        after the very last paragraph in the document. ---

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastRange.Collapse(0)
$lastRange.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Shading.BackgroundPatternColor = 16448250
